$wb = $excel.ActiveWorkbook

# --- LoginData sheet: add two new negative-login test case rows ---
$ws1 = $wb.Worksheets.Item("LoginData")

# Row 4: Login_Invalid_Username scenario (wrong username, correct password)
$ws1.Range("A3:D3").Copy($ws1.Range("A4:D4"))
$ws1.Range("A4").Value = "Login_Invalid_Username"
$ws1.Range("B4").Value = "wrongUser"
$ws1.Range("C4").Value = "Sdet@2025"
$ws1.Range("D4").Value = "Failure"

# Row 5: Login_Invalid_Password scenario (correct username, wrong password)
$ws1.Range("A3:D3").Copy($ws1.Range("A5:D5"))
$ws1.Range("A5").Value = "Login_Invalid_Password"
$ws1.Range("B5").Value = "Test_01"
$ws1.Range("C5").Value = "wrong123"
$ws1.Range("D5").Value = "Failure"

# --- EditorData sheet: rename RowIndex column to TestCaseName and switch
#     the row-identifying values from numeric indices to descriptive names ---
$ws2 = $wb.Worksheets.Item("EditorData")

$ws2.Range("A1").Value = "TestCaseName"
$ws2.Range("A2").Value = "PythonCode_Valid"
$ws2.Range("A3").Value = "PythonCode_Invalid"
